$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Row 3 now documents the underground temperature parameter instead of the
# obsolete house electricity demand peak parameter
$ws.Range("B3").Value = "p_undergroundTemperature_degC"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "degC"
$ws.Range("E3").ClearContents()

# Remove obsolete parameter rows 4 through 10 entirely
$ws.Range("A4:E10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update the active selection to match the saved view state
$ws.Range("B7").Select()
